$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adjust scene camera offset properties (CamOffestPos / CamOffestRot columns J/K).
# Row 2 = villageScene
$ws.Range("J2").Value = "0,8,7"
$ws.Range("K2").Value = "45,180"

# Row 6 = City / SelectScene
$ws.Range("J6").Value = "0,8,-7"
$ws.Range("K6").Value = "45,0"

# Match the authored selection state (cursor moved to K7 on save).
$ws.Range("K7").Select()
